$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 81.03634266666667
$ws.Range("H2").Value = 243.109028
$ws.Range("I2").Value = 0.1632931649012984
$ws.Range("J2").Value = 0.1632931649012984
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 15.35884066666667
$ws.Range("N2").Value = 46.076522
$ws.Range("O2").Value = 0.1012042817263867
$ws.Range("P2").Value = 0.1012042817263867
$ws.Range("Q2").Value = 1244.624275226735
$ws.Range("R2").Value = 11201.61847704062
$ws.Range("S2").Value = 0.01652596746466431
$ws.Range("T2").Value = 0.01652596746466431
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 81.03634266666667
$ws.Range("H3").Value = 243.109028
$ws.Range("I3").Value = 0.1632931649012984
$ws.Range("J3").Value = 0.1632931649012984
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 50.59256466666667
$ws.Range("N3").Value = 151.777694
$ws.Range("O3").Value = 0.3333704853712116
$ws.Range("P3").Value = 0.3333704853712116
$ws.Range("Q3").Value = 4099.836406713493
$ws.Range("R3").Value = 36898.52766042144
$ws.Range("S3").Value = 0.05443712164094712
$ws.Range("T3").Value = 0.05443712164094712
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 81.03634266666667
$ws.Range("H4").Value = 243.109028
$ws.Range("I4").Value = 0.1632931649012984
$ws.Range("J4").Value = 0.1632931649012984
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 60.37715666666667
$ws.Range("N4").Value = 181.13147
$ws.Range("O4").Value = 0.397844271305776
$ws.Range("P4").Value = 0.397844271305776
$ws.Range("Q4").Value = 4892.743956879019
$ws.Range("R4").Value = 44034.69561191116
$ws.Range("S4").Value = 0.06496525019937097
$ws.Range("T4").Value = 0.06496525019937097
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 81.03634266666667
$ws.Range("H5").Value = 243.109028
$ws.Range("I5").Value = 0.1632931649012984
$ws.Range("J5").Value = 0.1632931649012984
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.43221733333333
$ws.Range("N5").Value = 76.29665199999999
$ws.Range("O5").Value = 0.1675809615966257
$ws.Range("P5").Value = 0.1675809615966258
$ws.Range("Q5").Value = 2060.933878597139
$ws.Range("R5").Value = 18548.40490737426
$ws.Range("S5").Value = 0.02736482559631596
$ws.Range("T5").Value = 0.02736482559631596
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 178.365814
$ws.Range("H6").Value = 535.097442
$ws.Range("I6").Value = 0.3594179761796791
$ws.Range("J6").Value = 0.3594179761796791
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 15.35884066666667
$ws.Range("N6").Value = 46.076522
$ws.Range("O6").Value = 0.1012042817263867
$ws.Range("P6").Value = 0.1012042817263867
$ws.Range("Q6").Value = 2739.492117606303
$ws.Range("R6").Value = 24655.42905845672
$ws.Range("S6").Value = 0.03637463811881597
$ws.Range("T6").Value = 0.03637463811881597
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 178.365814
$ws.Range("H7").Value = 535.097442
$ws.Range("I7").Value = 0.3594179761796791
$ws.Range("J7").Value = 0.3594179761796791
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 50.59256466666667
$ws.Range("N7").Value = 151.777694
$ws.Range("O7").Value = 0.3333704853712116
$ws.Range("P7").Value = 0.3333704853712116
$ws.Range("Q7").Value = 9023.983979117638
$ws.Range("R7").Value = 81215.85581205874
$ws.Range("S7").Value = 0.1198193451701582
$ws.Range("T7").Value = 0.1198193451701582
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 178.365814
$ws.Range("H8").Value = 535.097442
$ws.Range("I8").Value = 0.3594179761796791
$ws.Range("J8").Value = 0.3594179761796791
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 60.37715666666667
$ws.Range("N8").Value = 181.13147
$ws.Range("O8").Value = 0.397844271305776
$ws.Range("P8").Value = 0.397844271305776
$ws.Range("Q8").Value = 10769.22069585553
$ws.Range("R8").Value = 96922.98626269975
$ws.Range("S8").Value = 0.1429923828274012
$ws.Range("T8").Value = 0.1429923828274012
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 178.365814
$ws.Range("H9").Value = 535.097442
$ws.Range("I9").Value = 0.3594179761796791
$ws.Range("J9").Value = 0.3594179761796791
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 25.43221733333333
$ws.Range("N9").Value = 76.29665199999999
$ws.Range("O9").Value = 0.1675809615966257
$ws.Range("P9").Value = 0.1675809615966258
$ws.Range("Q9").Value = 4536.238146484909
$ws.Range("R9").Value = 40826.14331836418
$ws.Range("S9").Value = 0.06023161006330375
$ws.Range("T9").Value = 0.06023161006330376
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 146.2303303333333
$ws.Range("H10").Value = 438.6909910000001
$ws.Range("I10").Value = 0.2946630198121519
$ws.Range("J10").Value = 0.2946630198121519
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 15.35884066666667
$ws.Range("N10").Value = 46.076522
$ws.Range("O10").Value = 0.1012042817263867
$ws.Range("P10").Value = 0.1012042817263867
$ws.Range("Q10").Value = 2245.9283442237
$ws.Range("R10").Value = 20213.3550980133
$ws.Range("S10").Value = 0.02982115927141687
$ws.Range("T10").Value = 0.02982115927141688
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 146.2303303333333
$ws.Range("H11").Value = 438.6909910000001
$ws.Range("I11").Value = 0.2946630198121519
$ws.Range("J11").Value = 0.2946630198121519
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 50.59256466666667
$ws.Range("N11").Value = 151.777694
$ws.Range("O11").Value = 0.3333704853712116
$ws.Range("P11").Value = 0.3333704853712116
$ws.Range("Q11").Value = 7398.167443617195
$ws.Range("R11").Value = 66583.50699255476
$ws.Range("S11").Value = 0.09823195393572401
$ws.Range("T11").Value = 0.09823195393572401
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 146.2303303333333
$ws.Range("H12").Value = 438.6909910000001
$ws.Range("I12").Value = 0.2946630198121519
$ws.Range("J12").Value = 0.2946630198121519
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 60.37715666666667
$ws.Range("N12").Value = 181.13147
$ws.Range("O12").Value = 0.397844271305776
$ws.Range("P12").Value = 0.397844271305776
$ws.Range("Q12").Value = 8828.971563954086
$ws.Range("R12").Value = 79460.74407558679
$ws.Range("S12").Value = 0.117229994397925
$ws.Range("T12").Value = 0.117229994397925
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 146.2303303333333
$ws.Range("H13").Value = 438.6909910000001
$ws.Range("I13").Value = 0.2946630198121519
$ws.Range("J13").Value = 0.2946630198121519
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 25.43221733333333
$ws.Range("N13").Value = 76.29665199999999
$ws.Range("O13").Value = 0.1675809615966257
$ws.Range("P13").Value = 0.1675809615966258
$ws.Range("Q13").Value = 3718.961541762459
$ws.Range("R13").Value = 33470.65387586213
$ws.Range("S13").Value = 0.04937991220708601
$ws.Range("T13").Value = 0.04937991220708601
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 90.630432
$ws.Range("H14").Value = 271.891296
$ws.Range("I14").Value = 0.1826258391068707
$ws.Range("J14").Value = 0.1826258391068707
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 15.35884066666667
$ws.Range("N14").Value = 46.076522
$ws.Range("O14").Value = 0.1012042817263867
$ws.Range("P14").Value = 0.1012042817263867
$ws.Range("Q14").Value = 1391.978364639168
$ws.Range("R14").Value = 12527.80528175251
$ws.Range("S14").Value = 0.0184825168714895
$ws.Range("T14").Value = 0.01848251687148951
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 90.630432
$ws.Range("H15").Value = 271.891296
$ws.Range("I15").Value = 0.1826258391068707
$ws.Range("J15").Value = 0.1826258391068707
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 50.59256466666667
$ws.Range("N15").Value = 151.777694
$ws.Range("O15").Value = 0.3333704853712116
$ws.Range("P15").Value = 0.3333704853712116
$ws.Range("Q15").Value = 4585.225991727936
$ws.Range("R15").Value = 41267.03392555143
$ws.Range("S15").Value = 0.06088206462438227
$ws.Range("T15").Value = 0.06088206462438227
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 90.630432
$ws.Range("H16").Value = 271.891296
$ws.Range("I16").Value = 0.1826258391068707
$ws.Range("J16").Value = 0.1826258391068707
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 60.37715666666667
$ws.Range("N16").Value = 181.13147
$ws.Range("O16").Value = 0.397844271305776
$ws.Range("P16").Value = 0.397844271305776
$ws.Range("Q16").Value = 5472.00779163168
$ws.Range("R16").Value = 49248.07012468512
$ws.Range("S16").Value = 0.07265664388107886
$ws.Range("T16").Value = 0.07265664388107886
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 90.630432
$ws.Range("H17").Value = 271.891296
$ws.Range("I17").Value = 0.1826258391068707
$ws.Range("J17").Value = 0.1826258391068707
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 25.43221733333333
$ws.Range("N17").Value = 76.29665199999999
$ws.Range("O17").Value = 0.1675809615966257
$ws.Range("P17").Value = 0.1675809615966258
$ws.Range("Q17").Value = 2304.932843637888
$ws.Range("R17").Value = 20744.39559274099
$ws.Range("S17").Value = 0.03060461372992005
$ws.Range("T17").Value = 0.03060461372992005

Write-Output "Applied 224 cell updates"
